$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.654.15'
$ws.Range("E2").Value = '  -4.13%  '

$ws.Range("D3").Value = '3.090.60'
$ws.Range("E3").Value = '  -5.47%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '605.49'
$ws.Range("E5").Value = '  -1.72%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.80'
$ws.Range("E6").Value = '  -8.94%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").Value = '3.085.09'
$ws.Range("E8").Value = '  -5.60%  '

$ws.Range("D9").Value = '0.518'
$ws.Range("E9").Value = '  -5.05%  '

$ws.Range("D10").Value = '0.149'
$ws.Range("E10").Value = '  -7.79%  '

$ws.Range("D11").Value = '5.21'
$ws.Range("E11").Value = '  -10.07%  '

$ws.Range("D12").Value = '0.466'
$ws.Range("E12").Value = '  -5.91%  '

$ws.Range("D13").Value = '0.0000247'
$ws.Range("E13").Value = '  -9.01%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.00'
$ws.Range("E14").Value = '  -10.38%  '

$ws.Range("D15").Value = '3.595.51'
$ws.Range("E15").Value = '  -5.43%  '

$ws.Range("D16").Value = '0.114'
$ws.Range("E16").Value = '  +0.71%  '

$ws.Range("D17").Value = '63.677.00'
$ws.Range("E17").Value = '  -4.21%  '

$ws.Range("D18").Value = '3.089.37'
$ws.Range("E18").Value = '  -5.42%  '

$ws.Range("D19").Value = '6.77'
$ws.Range("E19").Value = '  -9.12%  '

$ws.Range("D20").Value = '473.02'
$ws.Range("E20").Value = '  -6.33%  '

$ws.Range("D21").Value = '14.53'
$ws.Range("E21").Value = '  -6.07%  '

$ws.Range("E22").Value = '  -8.03%  '

$ws.Range("D23").Value = '7.64'
$ws.Range("E23").Value = '  -6.16%  '

$ws.Range("D24").Value = '13.47'
$ws.Range("E24").Value = '  -8.21%  '

$ws.Range("D25").Value = '83.11'
$ws.Range("E25").Value = '  -4.56%  '

$ws.Range("E26").Value = '  +0.07%  '

$ws.Range("E27").Value = '  -9.90%  '

$ws.Range("D28").Value = '8.32'
$ws.Range("E28").Value = '  -9.82%  '

$ws.Range("E29").Value = '  -11.62%  '

$ws.Range("D30").Value = '6.65'
$ws.Range("E30").Value = '  -5.31%  '

$ws.Range("D31").Value = '0.112'
$ws.Range("E31").Value = '  -12.25%  '

$ws.Range("E32").Value = '  +0.07%  '

$ws.Range("D33").Value = '2.69'
$ws.Range("E33").Value = '  -6.63%  '

$ws.Range("D34").Value = '26.03'
$ws.Range("E34").Value = '  -7.16%  '

$ws.Range("E35").Value = '  -4.55%  '

$ws.Range("D36").Value = '5.89'
$ws.Range("E36").Value = '  -9.19%  '

$ws.Range("D37").Value = '52.37'
$ws.Range("E37").Value = '  -5.92%  '

$ws.Range("D38").Value = '0.0₃0727'
$ws.Range("E38").Value = '  -8.20%  '

$ws.Range("D39").Value = '453.97'
$ws.Range("E39").Value = '  -8.22%  '

$ws.Range("D40").Value = '2.89'
$ws.Range("E40").Value = '  -15.06%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0390'
$ws.Range("E41").Value = '  -7.95%  '

$ws.Range("E42").Value = '  -8.35%  '

$ws.Range("D43").Value = '8.29'
$ws.Range("E43").Value = '  -6.19%  '

$ws.Range("D44").Value = '2.815.71'
$ws.Range("E44").Value = '  -6.59%  '

$ws.Range("D45").Value = '0.265'
$ws.Range("E45").Value = '  -9.81%  '

$ws.Range("D46").Value = '2.23'
$ws.Range("E46").Value = '  -12.16%  '

$ws.Range("E47").Value = '  -5.17%  '

$ws.Range("D49").Value = '25.94'
$ws.Range("E49").Value = '  -10.46%  '

$ws.Range("E50").Value = '  -5.85%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '117.60'
$ws.Range("E51").Value = '  -2.55%  '
